$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D width: widen from ~8.43 to ~10.43, drop "best fit" (explicit width now)
$ws.Columns.Item(4).ColumnWidth = 9.6

# A3 previously held "." - clear it out
$ws.Range("A3").Value = ""

# D3 formula now references row 2 (O2) instead of row 3 (O3)
$ws.Range("D3").Formula = "=IF(O2=0,0,IF(O2=21,1,3))"

# M3 / N3 now hold the literal text "date"
$ws.Range("M3").Value = "date"
$ws.Range("N3").Value = "date"
